$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1, copy style from G1 (bold header style)
$ws.Range("H1").Value = "Save"
$ws.Range("H1").Style = $ws.Range("G1").Style

# Save column values for rows 2-6
$values = @(1, 0, 1, 0, 1)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
